$d = $word.ActiveDocument

# Replace the full body content with the edited OOXML (title/author/email swap,
# rewritten essay body discussing government instead of arts/AI, new "Body:"
# section with elections/economy/international-relations paragraphs, rewritten
# summary, and a trailing empty paragraph) while leaving the final <w:sectPr/>
# (page setup) untouched.
$xml = @'
<w:p w:rsidR="00885AA9" w:rsidRDefault="001C74AC"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="44"/></w:rPr><w:t>Exploring the Dynamics of Government: A Path to Informed Citizenship</w:t></w:r></w:p><w:p w:rsidR="00885AA9" w:rsidRDefault="001C74AC"><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>Evelyn Mercer</w:t></w:r></w:p><w:p w:rsidR="00885AA9" w:rsidRDefault="001C74AC"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>mercerevelyn87@yahoo</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>com</w:t></w:r></w:p><w:p w:rsidR="00885AA9" w:rsidRDefault="00885AA9"/><w:p w:rsidR="00885AA9" w:rsidRDefault="001C74AC"><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>Throughout history, governments have been the framework that guides societies and influences the lives of individuals</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> As future citizens and leaders, understanding the mechanisms and functions of government is essential for navigating the complexities of our political landscape and fulfilling our civic duties</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> This exploration entails examining various components of government, analyzing their interrelationships, and comprehending the implications for citizens</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>Delving into the different forms and structures of governments, we can unravel the unique characteristics of each system</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Whether it's a direct democracy, where citizens directly make decisions, or a representative democracy, where elected officials act on behalf of the people, the design of a government significantly impacts its decision-making processes and outcomes</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Furthermore, navigating the intricate web of branches of government - legislative, executive, and judicial - enables us to decipher the balance of power among different governing bodies</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>Exploring the role of government in ensuring public well-being is a crucial aspect of this study</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Examining policies related to education, healthcare, infrastructure, and welfare unravels the ways in which governments allocate resources, regulate industries, and provide essential services to their citizens</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Understanding the delicate balance between individual liberty and collective responsibility allows us to critically evaluate the impact of government interventions on society</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>Body:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>Comprehending the dynamics of elections and political participation empowers citizens to actively shape their government</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Analyzing various electoral systems, such as first-past-the-post or proportional representation, illuminates how votes translate into representation</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Exploring the role of political parties, interest groups, and media in shaping public opinion reveals the dynamics that influence </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>government behavior</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Moreover, examining the rights and responsibilities of citizens - such as voting, holding public office, and engaging in civil discourse - emphasizes the importance of civic engagement in fostering accountable and responsive governance</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>Investigating the relationship between government and the economy provides insights into the intricate interplay of policies and economic outcomes</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Exploring concepts such as taxation, budgeting, and fiscal and monetary policies unveils how governments influence economic growth, employment, and inflation</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Analyzing the regulatory role of government in various sectors - such as finance, healthcare, and energy - helps understand the impact of regulations on business conduct and consumer welfare</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Delving into the complex relationship between government spending and economic development facilitates an exploration of how public expenditures impact infrastructure, education, and social welfare programs</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>The globalized world we live in demands an understanding of international relations and the role of government in this domain</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Examining international organizations, such as the United Nations, the World Bank, and the International Monetary Fund, uncovers how governments collaborate to address global challenges like climate change, poverty, and trade disputes</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Analyzing diplomatic relations, treaties, and alliances illuminates the strategies employed by governments to pursue their national interests in the international arena</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Moreover, exploring the evolving concept of globalization and its impact on national sovereignty unveils the complexities of interdependence and cooperation in the modern world</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p w:rsidR="00885AA9" w:rsidRDefault="001C74AC"><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="28"/></w:rPr><w:t>Summary</w:t></w:r></w:p><w:p w:rsidR="00885AA9" w:rsidRDefault="001C74AC"><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>In conclusion, understanding the multifaceted dynamics of government empowers individuals with the knowledge to navigate the complexities of their political system, exercise their civic duties effectively, and critically evaluate the role of government in shaping their lives and the world around them</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> From analyzing forms and structures of government to examining the delicate balance between individual liberty and collective responsibility, this exploration fosters informed citizenship and encourages active participation in the democratic process</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> By unraveling the intricate web of governmental mechanisms, citizens can navigate the political landscape with confidence and contribute to the well-being of their communities and the world</w:t></w:r><w:r w:rsidR="007827B0"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p/>
'@

$d.Content.InsertXML($xml)
